$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1712.6666
$ws.Range("I137").Value = 1277.7778
$ws.Range("K137").Value = 3833.3334
$ws.Range("M137").Value = -1283.3334
$ws.Range("H138").Value = 1955.4348
$ws.Range("J138").Value = 3000
$ws.Range("L138").Value = 9000
$ws.Range("N138").Value = -19280

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8335300.5
$ws.Range("I2").Value = 12501026
$ws.Range("K2").Value = 12501026
$ws.Range("M2").Value = -12500913
$ws.Range("H45").Value = 2833.1667
$ws.Range("J45").Value = 4133.3335
$ws.Range("L45").Value = 4133.3335
$ws.Range("N45").Value = -4887.3335
$ws.Range("H76").Value = 19999
$ws.Range("J76").Value = 19999
$ws.Range("L76").Value = 19999
$ws.Range("N76").Value = -20675
$ws.Range("H79").Value = 19999
$ws.Range("J79").Value = 19999
$ws.Range("L79").Value = 19999
$ws.Range("N79").Value = -22339
$ws.Range("H116").Value = 8335300.5
$ws.Range("I116").Value = 12501026
$ws.Range("K116").Value = 12501026
$ws.Range("M116").Value = -12498732

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8335300.5
$ws.Range("I3").Value = 12501026
$ws.Range("K3").Value = 12501026
$ws.Range("M3").Value = -12500912
$ws.Range("H20").Value = 3413.4614
$ws.Range("I20").Value = 2204.7
$ws.Range("K20").Value = 2204.7
$ws.Range("M20").Value = -1957.7
$ws.Range("H29").Value = 2999.5
$ws.Range("J29").Value = 2999.5
$ws.Range("L29").Value = 2999.5
$ws.Range("N29").Value = -3577.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 18001
$ws.Range("I4").Value = 6000
$ws.Range("J4").Value = 30002
$ws.Range("K4").Value = 6000
$ws.Range("L4").Value = 30002
$ws.Range("M4").Value = -5888
$ws.Range("N4").Value = -30226
$ws.Range("H7").Value = 155.5625
$ws.Range("I7").Value = 128.90909
$ws.Range("J7").Value = 214.2
$ws.Range("K7").Value = 128.90909
$ws.Range("L7").Value = 214.2
$ws.Range("M7").Value = -15.90908999999999
$ws.Range("N7").Value = -440.2
$ws.Range("H15").Value = 3915
$ws.Range("J15").Value = 413.33334
$ws.Range("L15").Value = 413.33334
$ws.Range("N15").Value = -753.33334
$ws.Range("H22").Value = 91788.586
$ws.Range("I22").Value = 91042
$ws.Range("K22").Value = 91042
$ws.Range("M22").Value = -90692
$ws.Range("H31").Value = 3919.5518
$ws.Range("J31").Value = 6186.1113
$ws.Range("L31").Value = 6186.1113
$ws.Range("N31").Value = -6776.1113
$ws.Range("H34").Value = 3919.5518
$ws.Range("J34").Value = 6186.1113
$ws.Range("L34").Value = 6186.1113
$ws.Range("N34").Value = -6590.1113
$ws.Range("H86").Value = 3225
$ws.Range("I86").Value = 3200
$ws.Range("J86").Value = 3250
$ws.Range("K86").Value = 3200
$ws.Range("L86").Value = 3250
$ws.Range("M86").Value = -2077
$ws.Range("N86").Value = -5496
$ws.Range("H89").Value = 3225
$ws.Range("I89").Value = 3200
$ws.Range("J89").Value = 3250
$ws.Range("K89").Value = 16000
$ws.Range("L89").Value = 16250
$ws.Range("M89").Value = -10384
$ws.Range("N89").Value = -27482
$ws.Range("H99").Value = 12526.192
$ws.Range("I99").Value = 7836.2856
$ws.Range("J99").Value = 17997.75
$ws.Range("K99").Value = 7836.2856
$ws.Range("L99").Value = 17997.75
$ws.Range("M99").Value = -6338.2856
$ws.Range("N99").Value = -20993.75
$ws.Range("H126").Value = 12526.192
$ws.Range("I126").Value = 7836.2856
$ws.Range("J126").Value = 17997.75
$ws.Range("K126").Value = 23508.8568
$ws.Range("L126").Value = 53993.25
$ws.Range("M126").Value = -21038.8568
$ws.Range("N126").Value = -58933.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1612.5
$ws.Range("I11").Value = 1965.8572
$ws.Range("J11").Value = 1117.8
$ws.Range("K11").Value = 5897.571599999999
$ws.Range("L11").Value = 3353.4
$ws.Range("M11").Value = -5757.571599999999
$ws.Range("N11").Value = -3633.4
$ws.Range("H23").Value = 200144.8
$ws.Range("J23").Value = 200144.8
$ws.Range("L23").Value = 600434.3999999999
$ws.Range("N23").Value = -600904.3999999999
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H109").Value = 1482
$ws.Range("I109").Value = 875
$ws.Range("J109").Value = 2999.5
$ws.Range("K109").Value = 2625
$ws.Range("L109").Value = 8998.5
$ws.Range("M109").Value = -1585
$ws.Range("N109").Value = -11078.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 10000583
$ws.Range("J3").Value = 874.5
$ws.Range("L3").Value = 874.5
$ws.Range("N3").Value = -1106.5
$ws.Range("H13").Value = 124.5
$ws.Range("I13").Value = 124.5
$ws.Range("K13").Value = 124.5
$ws.Range("M13").Value = 14.5
$ws.Range("H55").Value = 5750
$ws.Range("J55").Value = 8000
$ws.Range("L55").Value = 8000
$ws.Range("N55").Value = -8654
$ws.Range("H70").Value = 5422.5557
$ws.Range("I70").Value = 5060.6
$ws.Range("J70").Value = 5875
$ws.Range("K70").Value = 5060.6
$ws.Range("L70").Value = 5875
$ws.Range("M70").Value = -4790.6
$ws.Range("N70").Value = -6415
$ws.Range("H73").Value = 5422.5557
$ws.Range("I73").Value = 5060.6
$ws.Range("J73").Value = 5875
$ws.Range("K73").Value = 5060.6
$ws.Range("L73").Value = 5875
$ws.Range("M73").Value = -4124.6
$ws.Range("N73").Value = -7747
$ws.Range("H93").Value = 49993.332
$ws.Range("J93").Value = 49993.332
$ws.Range("L93").Value = 49993.332
$ws.Range("N93").Value = -53737.332
$ws.Range("H97").Value = 2285.3572
$ws.Range("I97").Value = 1368
$ws.Range("J97").Value = 2973.375
$ws.Range("K97").Value = 1368
$ws.Range("L97").Value = 2973.375
$ws.Range("M97").Value = -872
$ws.Range("N97").Value = -3965.375

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 3908.3333
$ws.Range("I4").Value = 3859
$ws.Range("K4").Value = 3859
$ws.Range("M4").Value = -3746
$ws.Range("H28").Value = 3908.3333
$ws.Range("I28").Value = 3859
$ws.Range("K28").Value = 3859
$ws.Range("M28").Value = -3627
$ws.Range("H37").Value = 3908.3333
$ws.Range("I37").Value = 3859
$ws.Range("K37").Value = 3859
$ws.Range("M37").Value = -3752
$ws.Range("H55").Value = 419.5
$ws.Range("J55").Value = 446.8
$ws.Range("L55").Value = 446.8
$ws.Range("N55").Value = -792.8
$ws.Range("H61").Value = 3958.8
$ws.Range("I61").Value = 3958.8
$ws.Range("K61").Value = 3958.8
$ws.Range("M61").Value = -3756.8
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H93").Value = 1515.5
$ws.Range("I93").Value = 898.25
$ws.Range("J93").Value = 2750
$ws.Range("K93").Value = 898.25
$ws.Range("L93").Value = 2750
$ws.Range("M93").Value = 349.75
$ws.Range("N93").Value = -5246
$ws.Range("H100").Value = 5699.2856
$ws.Range("I100").Value = 5315.8335
$ws.Range("K100").Value = 5315.8335
$ws.Range("M100").Value = -4774.8335
$ws.Range("H113").Value = 3958.8
$ws.Range("I113").Value = 3958.8
$ws.Range("K113").Value = 3958.8
$ws.Range("M113").Value = -1788.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 10002250
$ws.Range("I10").Value = 20000000
$ws.Range("J10").Value = 4500
$ws.Range("K10").Value = 20000000
$ws.Range("L10").Value = 4500
$ws.Range("M10").Value = -19999831
$ws.Range("N10").Value = -4838
$ws.Range("H12").Value = 8599
$ws.Range("I12").Value = 2898.5
$ws.Range("K12").Value = 2898.5
$ws.Range("M12").Value = -2756.5
$ws.Range("H21").Value = 40017
$ws.Range("J21").Value = 40017
$ws.Range("L21").Value = 40017
$ws.Range("N21").Value = -40487
$ws.Range("H35").Value = 40017
$ws.Range("J35").Value = 40017
$ws.Range("L35").Value = 40017
$ws.Range("N35").Value = -40597
$ws.Range("H40").Value = 32129.5
$ws.Range("J40").Value = 32129.5
$ws.Range("L40").Value = 32129.5
$ws.Range("N40").Value = -32427.5
$ws.Range("H45").Value = 17626
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 17626
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 17626
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -18608
$ws.Range("H100").Value = 2029.7142
$ws.Range("I100").Value = 2205.5715
$ws.Range("J100").Value = 1678
$ws.Range("K100").Value = 4411.143
$ws.Range("L100").Value = 3356
$ws.Range("M100").Value = -3870.143
$ws.Range("N100").Value = -4438
$ws.Range("H113").Value = 3844.8333
$ws.Range("I113").Value = 767.5
$ws.Range("J113").Value = 9999.5
$ws.Range("K113").Value = 2302.5
$ws.Range("L113").Value = 29998.5
$ws.Range("M113").Value = -132.5
$ws.Range("N113").Value = -34338.5
$ws.Range("H136").Value = 2522.3076
$ws.Range("I136").Value = 1415.125
$ws.Range("K136").Value = 4245.375
$ws.Range("M136").Value = -1695.375
